# Update the cryptos price/volume table with the latest scraped figures
# (GitHub Actions scheduled refresh). A handful of "Price" cells look like
# plain numbers (e.g. "517.12") but must stay as text, exactly as the
# source data stores them, so a leading apostrophe (quote-prefix) is used
# for those to stop Excel from auto-converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.099.08'
$ws.Range('E2').Value = '  +3.62%  '
$ws.Range('D3').Value = '2.317.48'
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''517.12'
$ws.Range('E5').Value = '  +4.25%  '
$ws.Range('D6').Value = '''133.27'
$ws.Range('E6').Value = '  +3.66%  '
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').Value = '''0.535'
$ws.Range('E8').Value = '  +1.98%  '
$ws.Range('D9').Value = '2.342.51'
$ws.Range('E9').Value = '  +3.43%  '
$ws.Range('E10').Value = '  +8.64%  '
$ws.Range('D11').Value = '''0.155'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').Value = '''5.16'
$ws.Range('E12').Value = '  +8.19%  '
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('D14').Value = '''24.07'
$ws.Range('E14').Value = '  +5.81%  '
$ws.Range('D15').Value = '2.731.02'
$ws.Range('E15').Value = '  +2.66%  '
$ws.Range('D16').Value = '56.179.60'
$ws.Range('E16').Value = '  +3.82%  '
$ws.Range('E17').Value = '  +4.60%  '
$ws.Range('D18').Value = '2.335.63'
$ws.Range('E18').Value = '  +3.19%  '
$ws.Range('D19').Value = '''10.54'
$ws.Range('E19').Value = '  +2.93%  '
$ws.Range('E20').Value = '  +2.86%  '
$ws.Range('D21').Value = '''321.78'
$ws.Range('E21').Value = '  +6.58%  '
$ws.Range('D22').Value = '''6.65'
$ws.Range('E22').Value = '  +5.55%  '
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '''60.67'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').Value = '''0.991'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('E26').Value = '  +6.32%  '
$ws.Range('E27').Value = '  +4.90%  '
$ws.Range('D28').Value = '''172.19'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').Value = '''1.19'
$ws.Range('E29').Value = '  +9.83%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0724'
$ws.Range('E30').Value = '  +4.96%  '
$ws.Range('E31').Value = '  +5.55%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '''1.68'
$ws.Range('E32').Value = '  +4.72%  '
$ws.Range('D33').Value = '''18.37'
$ws.Range('E33').Value = '  +3.73%  '
$ws.Range('D34').Value = '''0.999'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E36').Value = '  +5.93%  '
$ws.Range('D37').Value = '''0.926'
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('E38').Value = '  +7.30%  '
$ws.Range('E39').Value = '  +8.76%  '
$ws.Range('D40').Value = '''37.37'
$ws.Range('E40').Value = '  +4.05%  '
$ws.Range('D41').Value = '''0.383'
$ws.Range('E41').Value = '  +2.96%  '
$ws.Range('D42').Value = '''3.63'
$ws.Range('E42').Value = '  +8.05%  '
$ws.Range('D43').Value = '''137.93'
$ws.Range('E43').Value = '  +10.04%  '
$ws.Range('E44').Value = '  +7.14%  '
$ws.Range('D45').Value = '''267.66'
$ws.Range('E45').Value = '  +10.99%  '
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('D47').Value = '''0.0926'
$ws.Range('E47').Value = '  +4.05%  '
$ws.Range('D48').Value = '''0.557'
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('D49').Value = '''0.382'
$ws.Range('E49').Value = '  +2.84%  '
$ws.Range('E50').Value = '  +5.69%  '
$ws.Range('D51').Value = '''16.93'
$ws.Range('E51').Value = '  +5.14%  '
